$d = $word.ActiveDocument

# --- Step 1: build the new "Play Fruits Free..." (bold) paragraph near the end
#     of the document, right after the "What we don't like" bullet list and
#     before the closing "Create a feature image..." paragraph. We clone the
#     existing "Meta description" paragraph's formatted content (so the new
#     paragraph picks up the same leading empty-run + bold-run layout) and
#     then re-point its text at the title. ---

$metaPara = $d.Paragraphs.Item(2)

$n = $d.Paragraphs.Count
$lastBulletPara = $d.Paragraphs.Item($n - 1)   # "Minimal bonus features outside of Lightning Rounds"

$lastBulletPara.Range.InsertParagraphAfter() | Out-Null

$n2 = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($n2 - 1)
$newPara.Style = "Normal"

$srcRange = $d.Range($metaPara.Range.Start, $metaPara.Range.End - 1)
$destRange = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$destRange.FormattedText = $srcRange.FormattedText

$newParaText = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$newParaText.Text = "Play Fruits Free - Review of Nolimit City's 6-Reel Slot"
$newParaBold = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$newParaBold.Font.Bold = 1

# --- Step 2: remove the original "Meta description" paragraph near the top. ---

$metaPara.Range.Delete()

# --- Step 3: the old final paragraph ("Create a feature image...") becomes the
#     meta-description text, keeping its italic run formatting. ---

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastParaText = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
$lastParaText.Text = "Try Fruits slot by Nolimit City with adaptable reels, high payout potential, and the unique Lightning Rounds feature. Play for free with our review."

Write-Output "done"
